# Weekly fruit/vegetable price update: a new daily price record is inserted
# into the dataset at row 421 (pushing the existing rows 421-492 down to
# 422-493), growing the used range from A1:R492 to A1:R493.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 421, shifting rows 421..492 down to 422..493
# (this is what naturally happened in the source system when the new
# weekly record was appended/sorted into this position).
$ws.Rows(421).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A421").Value = 5
$ws.Range("B421").Value = "Macroferia Regional de Talca"
$ws.Range("C421").Value = "Maule"
$ws.Range("D421").Value = 44951
$ws.Range("E421").Value = 7
$ws.Range("F421").Value = 100112032
$ws.Range("G421").Value = "Zapallo italiano"
$ws.Range("H421").Value = "Sin especificar"
$ws.Range("I421").Value = "Primera"
$ws.Range("J421").Value = 400
$ws.Range("K421").Value = 5000
$ws.Range("L421").Value = 5000
$ws.Range("M421").Value = 5000
$ws.Range("N421").Value = "$/caja 50 unidades"
$ws.Range("O421").Value = "Región del Maule"
$ws.Range("P421").Value = 100
$ws.Range("Q421").Value = 50
$ws.Range("R421").Value = "Hortaliza"
